$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 733.3333
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H23").Value = 733.3333
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H29").Value = 1438
$ws.Range("I29").Value = 297.5
$ws.Range("K29").Value = 892.5
$ws.Range("M29").Value = -611.5

$ws.Range("H58").Value = 2235.2856
$ws.Range("J58").Value = 6900
$ws.Range("L58").Value = 20700
$ws.Range("N58").Value = -21000

$ws.Range("H86").Value = 14326963
$ws.Range("I86").Value = 7872.5
$ws.Range("J86").Value = 20054600
$ws.Range("K86").Value = 7872.5
$ws.Range("L86").Value = 20054600
$ws.Range("M86").Value = -6749.5
$ws.Range("N86").Value = -20056846

$ws.Range("H89").Value = 14326963
$ws.Range("I89").Value = 7872.5
$ws.Range("J89").Value = 20054600
$ws.Range("K89").Value = 39362.5
$ws.Range("L89").Value = 100273000
$ws.Range("M89").Value = -33746.5
$ws.Range("N89").Value = -100284232

$ws.Range("H111").Value = 41241.5
$ws.Range("I111").Value = 18973.8
$ws.Range("K111").Value = 56921.39999999999
$ws.Range("M111").Value = -53854.39999999999

$ws.Range("H112").Value = 126763.375
$ws.Range("J112").Value = 2018.1666
$ws.Range("L112").Value = 6054.4998
$ws.Range("N112").Value = -8270.4998

$ws.Range("H137").Value = 3750.3
$ws.Range("I137").Value = 2563.5
$ws.Range("J137").Value = 4181.864
$ws.Range("K137").Value = 7690.5
$ws.Range("L137").Value = 12545.592
$ws.Range("M137").Value = -5140.5
$ws.Range("N137").Value = -17645.592

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 373669.12
$ws.Range("I32").Value = 420090.28
$ws.Range("K32").Value = 420090.28
$ws.Range("M32").Value = -419803.28

$ws.Range("H61").Value = 1776.5454
$ws.Range("I61").Value = 1504
$ws.Range("K61").Value = 1504
$ws.Range("M61").Value = -1292

$ws.Range("H74").Value = 2641.3408
$ws.Range("J74").Value = 2410.875
$ws.Range("L74").Value = 2410.875
$ws.Range("N74").Value = -4158.875

$ws.Range("H77").Value = 2641.3408
$ws.Range("J77").Value = 2410.875
$ws.Range("L77").Value = 12054.375
$ws.Range("N77").Value = -20790.375

$ws.Range("H110").Value = 90924250
$ws.Range("I110").Value = 142858670
$ws.Range("K110").Value = 142858670
$ws.Range("M110").Value = -142856625

$ws.Range("H122").Value = 19611244
$ws.Range("I122").Value = 33335516
$ws.Range("K122").Value = 100006548
$ws.Range("M122").Value = -100004098

$ws.Range("H136").Value = 1776.5454
$ws.Range("I136").Value = 1504
$ws.Range("K136").Value = 4512
$ws.Range("M136").Value = -1962

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 5757.25
$ws.Range("I37").Value = 1000
$ws.Range("K37").Value = 1000
$ws.Range("M37").Value = -863

$ws.Range("H86").Value = 1905.4445
$ws.Range("I86").Value = 1760.1875
$ws.Range("K86").Value = 1760.1875
$ws.Range("M86").Value = -637.1875

$ws.Range("H89").Value = 1905.4445
$ws.Range("I89").Value = 1760.1875
$ws.Range("K89").Value = 8800.9375
$ws.Range("M89").Value = -3184.9375

$ws.Range("H107").Value = 15164432
$ws.Range("J107").Value = 71459860
$ws.Range("L107").Value = 71459860
$ws.Range("N107").Value = -71463700

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3409
$ws.Range("I31").Value = 2568
$ws.Range("J31").Value = 3813.926
$ws.Range("K31").Value = 2568
$ws.Range("L31").Value = 3813.926
$ws.Range("M31").Value = -2273
$ws.Range("N31").Value = -4403.925999999999

$ws.Range("H34").Value = 3409
$ws.Range("I34").Value = 2568
$ws.Range("J34").Value = 3813.926
$ws.Range("K34").Value = 2568
$ws.Range("L34").Value = 3813.926
$ws.Range("M34").Value = -2366
$ws.Range("N34").Value = -4217.925999999999

$ws.Range("H122").Value = 1924.2354
$ws.Range("I122").Value = 1274.2
$ws.Range("K122").Value = 3822.6
$ws.Range("M122").Value = -1372.6

$ws.Range("H132").Value = 6464.5
$ws.Range("I132").Value = 2994.5
$ws.Range("K132").Value = 8983.5
$ws.Range("M132").Value = -6453.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 45014092
$ws.Range("I4").Value = 51938944
$ws.Range("K4").Value = 155816832
$ws.Range("M4").Value = -155816720

$ws.Range("H11").Value = 181234.67
$ws.Range("I11").Value = 61386.35
$ws.Range("K11").Value = 184159.05
$ws.Range("M11").Value = -184019.05

$ws.Range("H26").Value = 316
$ws.Range("J26").Value = 598
$ws.Range("L26").Value = 1794
$ws.Range("N26").Value = -2370

$ws.Range("H60").Value = 499.5
$ws.Range("I60").Value = 499.5
$ws.Range("K60").Value = 1498.5
$ws.Range("M60").Value = -1247.5

$ws.Range("H63").Value = 15996.8
$ws.Range("I63").Value = 16666
$ws.Range("J63").Value = 14993
$ws.Range("K63").Value = 49998
$ws.Range("L63").Value = 44979
$ws.Range("M63").Value = -49249
$ws.Range("N63").Value = -46477

$ws.Range("H66").Value = 15996.8
$ws.Range("I66").Value = 16666
$ws.Range("J66").Value = 14993
$ws.Range("K66").Value = 149994
$ws.Range("L66").Value = 134937
$ws.Range("M66").Value = -146250
$ws.Range("N66").Value = -142425

$ws.Range("H68").Value = 1896.4375
$ws.Range("J68").Value = 2374.625
$ws.Range("L68").Value = 7123.875
$ws.Range("N68").Value = -8745.875

$ws.Range("H71").Value = 1896.4375
$ws.Range("J71").Value = 2374.625
$ws.Range("L71").Value = 21371.625
$ws.Range("N71").Value = -29483.625

$ws.Range("H129").Value = 2030.4762
$ws.Range("J129").Value = 2575
$ws.Range("L129").Value = 7725
$ws.Range("N129").Value = -17725

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3293.5833
$ws.Range("I80").Value = 1997.5
$ws.Range("K80").Value = 1997.5
$ws.Range("M80").Value = -999.5

$ws.Range("H83").Value = 3293.5833
$ws.Range("I83").Value = 1997.5
$ws.Range("K83").Value = 9987.5
$ws.Range("M83").Value = -4995.5

$ws.Range("H102").Value = 2314.2593
$ws.Range("I102").Value = 1444.3
$ws.Range("K102").Value = 1444.3
$ws.Range("M102").Value = 177.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1348.68
$ws.Range("J22").Value = 1357
$ws.Range("L22").Value = 1357
$ws.Range("N22").Value = -1947

$ws.Range("H27").Value = 1348.68
$ws.Range("J27").Value = 1357
$ws.Range("L27").Value = 1357
$ws.Range("N27").Value = -1571

$ws.Range("H82").Value = 2459.8125
$ws.Range("I82").Value = 1335.75
$ws.Range("J82").Value = 3583.875
$ws.Range("K82").Value = 1335.75
$ws.Range("L82").Value = 3583.875
$ws.Range("M82").Value = -974.75
$ws.Range("N82").Value = -4305.875

$ws.Range("H85").Value = 2459.8125
$ws.Range("I85").Value = 1335.75
$ws.Range("J85").Value = 3583.875
$ws.Range("K85").Value = 1335.75
$ws.Range("L85").Value = 3583.875
$ws.Range("M85").Value = -87.75
$ws.Range("N85").Value = -6079.875

$ws.Range("H136").Value = 5889.3228
$ws.Range("I136").Value = 3942.6667
$ws.Range("J136").Value = 9977.299999999999
$ws.Range("K136").Value = 11828.0001
$ws.Range("L136").Value = 29931.9
$ws.Range("M136").Value = -9278.000100000001
$ws.Range("N136").Value = -35031.89999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 464.18182
$ws.Range("I100").Value = 444.7143
$ws.Range("K100").Value = 889.4286
$ws.Range("M100").Value = -348.4286

$ws.Range("H113").Value = 3086.25
$ws.Range("I113").Value = 4476
$ws.Range("J113").Value = 1696.5
$ws.Range("K113").Value = 13428
$ws.Range("L113").Value = 5089.5
$ws.Range("M113").Value = -11258
$ws.Range("N113").Value = -9429.5
